# Script: apply scraper update to czech-republic cfl-group-a 2023-2024 sheet
# - Swaps the F:V (match detail) content between several pairs of rows whose
#   home/away teams had been attributed to the wrong row by the scraper.
# - Appends five newly scraped matches at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap F:V between paired rows (same match day / tournament, rows were
#    simply out of order relative to the canonical scrape order).
# ---------------------------------------------------------------------------
$pairs = @(
    @(2,3),
    @(14,15),
    @(16,17),
    @(30,31),
    @(70,71),
    @(79,80),
    @(94,95),
    @(102,103)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $rng1 = $ws.Range("F$r1`:V$r1")
    $rng2 = $ws.Range("F$r2`:V$r2")
    $val1 = $rng1.Value()
    $val2 = $rng2.Value()
    $rng1.Value = $val2
    $rng2.Value = $val1
}

# ---------------------------------------------------------------------------
# 2) Append five new rows (106-110), copying the formatting of the last
#    existing data row (105) and filling in the scraped values.
# ---------------------------------------------------------------------------
$lastRow = 105
$newRowsCount = 5

for ($i = 1; $i -le $newRowsCount; $i++) {
    $target = $lastRow + $i
    $src = $ws.Range("A$lastRow`:V$lastRow")
    $dst = $ws.Range("A$target`:V$target")
    $src.Copy($dst)
}

$newData = @(
    @{ Row=106; Indice=105; Data=45234.42708333334; Home="Hostoun"; HomeG=0; Away="Domazlice"; AwayG=3;
       J=2.3;  K="04/11/2023 00:42"; L=2.53; M="04/11/2023 10:10";
       N=3.43; O="04/11/2023 00:42"; P=3.59; Q="04/11/2023 10:10";
       R=2.66; S="04/11/2023 00:42"; T=2.43; U="04/11/2023 10:10";
       V="https://www.betexplorer.com/football/czech-republic/cfl-group-a/hostoun-domazlice/I3iVN3iI/" },

    @{ Row=107; Indice=106; Data=45234.42708333334; Home="Vltavin"; HomeG=2; Away="Bohemians 1905 B"; AwayG=2;
       J=1.97; K="04/11/2023 00:42"; L=2.03; M="04/11/2023 10:14";
       N=3.59; O="04/11/2023 00:42"; P=3.37; Q="04/11/2023 10:14";
       R=3.16; S="04/11/2023 00:42"; T=3.37; U="04/11/2023 10:14";
       V="https://www.betexplorer.com/football/czech-republic/cfl-group-a/loko-vltavin-bohemians-1905/hK8sRZIu/" },

    @{ Row=108; Indice=107; Data=45234.4375; Home="Kraluv Dvur"; HomeG=1; Away="Pisek"; AwayG=5;
       J=1.97; K="04/11/2023 01:13"; L=1.98; M="04/11/2023 10:27";
       N=3.59; O="04/11/2023 01:13"; P=3.53; Q="04/11/2023 10:27";
       R=3.16; S="04/11/2023 01:13"; T=3.37; U="04/11/2023 10:27";
       V="https://www.betexplorer.com/football/czech-republic/cfl-group-a/kraluv-dvur-pisek/CQjZMN6O/" },

    @{ Row=109; Indice=108; Data=45234.4375; Home="Ceske Budejovice B"; HomeG=0; Away="Slavia Prague B"; AwayG=0;
       J=2.89; K="04/11/2023 01:13"; L=3.28; M="04/11/2023 09:37";
       N=3.8;  O="04/11/2023 01:13"; P=4.13; Q="04/11/2023 10:25";
       R=2.03; S="04/11/2023 01:13"; T=1.85; U="04/11/2023 10:25";
       V="https://www.betexplorer.com/football/czech-republic/cfl-group-a/ceske-budejovice-slavia-prague/WA7oQgYo/" },

    @{ Row=110; Indice=109; Data=45234.58333333334; Home="FK Robstav"; HomeG=2; Away="Dukla Prague B"; AwayG=2;
       J=1.94; K="04/11/2023 12:13"; L=1.54; M="04/11/2023 13:52";
       N=3.65; O="04/11/2023 12:13"; P=4.25; Q="04/11/2023 13:52";
       R=3.25; S="04/11/2023 12:13"; T=5.01; U="04/11/2023 13:52";
       V="https://www.betexplorer.com/football/czech-republic/cfl-group-a/fk-robstav-dukla-prague/QTfwMsMU/" }
)

foreach ($row in $newData) {
    $r = $row.Row
    $ws.Range("A$r").Value = $row.Indice
    $ws.Range("B$r").Value = "czech-republic"
    $ws.Range("C$r").Value = "cfl-group-a"
    $ws.Range("D$r").Value = "2023-2024"
    $ws.Range("E$r").Value = $row.Data
    $ws.Range("F$r").Value = $row.Home
    $ws.Range("G$r").Value = $row.HomeG
    $ws.Range("H$r").Value = $row.Away
    $ws.Range("I$r").Value = $row.AwayG
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = $row.K
    $ws.Range("L$r").Value = $row.L
    $ws.Range("M$r").Value = $row.M
    $ws.Range("N$r").Value = $row.N
    $ws.Range("O$r").Value = $row.O
    $ws.Range("P$r").Value = $row.P
    $ws.Range("Q$r").Value = $row.Q
    $ws.Range("R$r").Value = $row.R
    $ws.Range("S$r").Value = $row.S
    $ws.Range("T$r").Value = $row.T
    $ws.Range("U$r").Value = $row.U
    $ws.Range("V$r").Value = $row.V
}
